$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("stickers")
$ws2 = $wb.Worksheets.Item("users")

# --- stickers sheet: drop the duplicate rows 8-10, keep only row 7 ---
$ws1.Rows.Item(8).Delete()
$ws1.Rows.Item(8).Delete()
$ws1.Rows.Item(8).Delete()

# --- users sheet: append the newly registered user ---
$ws2.Range("A2").Value = 1990838212
$ws2.Range("B2").Value = "Jack"
$ws2.Range("C2").Value = "м"
$ws2.Range("D2").Value = "8а"

# Update the remaining row's answer text (closing paren added)
$ws1.Range("C7").Value = "жизнь - интересная штука)"

# --- selection / active sheet bookkeeping ---
$null = $ws1.Range("C7").Select()
$null = $ws2.Range("A5:D5").Select()
$null = $ws1.Activate()
